$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before row 7 (they inherit some formatting from row 6,
# which is fixed up below) to make room for the new "FERRITIN_*" variables in
# the "Donor" file section. Everything at/after row 7 shifts down by 3.
$ws.Rows("7:9").Insert()

# Columns C/E of the freshly inserted rows come out unstyled, but column D
# inherits row 6's style (s=2). Re-align D7:D8 with the plain style (s=0)
# used by the rest of the "Donor" block by pasting the format from a cell
# that already uses that plain style; row 9 is meant to keep s=2, so it is
# left untouched.
$ws.Range("C3").Copy()
$ws.Range("D7:D8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 7: FERRITIN_FIRST
$ws.Range("B7").Value = "Donor"
$ws.Range("C7").Value = "FERRITIN_FIRST"
$ws.Range("D7").Value = "double"
$ws.Range("E7").Value = "First ferritin (ug/L = ng/mL)"

# Row 8: FERRITIN_LAST
$ws.Range("B8").Value = "Donor"
$ws.Range("C8").Value = "FERRITIN_LAST"
$ws.Range("D8").Value = "double"
$ws.Range("E8").Value = "Last ferritin (ug/L = ng/mL)"

# Row 9: FERRITIN_LAST_DATE
$ws.Range("B9").Value = "Donor"
$ws.Range("C9").Value = "FERRITIN_LAST_DATE"
$ws.Range("D9").Value = "varchar(n)"
$ws.Range("E9").Value = "Date when the last ferritin was measured (yyyymmdd)"

# Give E8 its own distinct cell style (matches the new style introduced by
# the original edit for that one cell).
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").Font.Size = 10

# Update the description of DONAT_RESULT_CODE (now on row 16 after the
# insert) to reflect the corrected unit (g/dL instead of g/L).
$ws.Range("E16").Value = "Hemoglogin (g/dL), for example 14.5"

# Leave the selection where the author left it when saving.
$ws.Range("E17").Select()
